# Update the estimation file: bump "Actual time" hours for the last two
# tasks (row 32: testing and bug fixing, row 33: cart-modal), and move the
# active selection so the sheet scrolls back to the top and the cursor
# rests on I23.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update actual-time values
$ws.Range("E32").Value = 6
$ws.Range("E33").Value = 3

# Update the current selection/view: selecting I23 moves the active cell
# there and resets the window so A1 (not A10) is the top-left cell again.
$ws.Range("I23").Select()
